{"js": "// Fill in the approval info in the \"Auteurs / Approbateurs / Validation\"\n// table near the top of the document:\n//   - the \"Approbateurs\" value cell (currently just a single space) gets\n//     the approver's name \"Fanny LAJEUNESSE\" appended;\n//   - the \"Approuv\u00e9 le :  \" value cell gets the approval date \"29/10\"\n//     appended.\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nif (tables.items.length === 0) {\n  throw new Error(\"Expected at least one table in the document body.\");\n}\n\nfor (const t of tables.items) {\n  t.load(\"values\");\n}\nawait context.sync();\n\n// Locate the (first, in document order) table that has an \"Approbateurs\"\n// header cell, and remember where, so the edit is resilient to the table\n// not being exactly where we expect it.\nlet infoTable = null;\nlet headerRowIndex = -1;\nlet approbateursColIndex = -1;\n\nfor (const t of tables.items) {\n  const values = t.values;\n  for (let r = 0; r < values.length && infoTable === null; r++) {\n    for (let c = 0; c < values[r].length; c++) {\n      const cellText = (values[r][c] || \"\").trim();\n      if (cellText === \"Approbateurs\") {\n        infoTable = t;\n        headerRowIndex = r;\n        approbateursColIndex = c;\n        break;\n      }\n    }\n  }\n  if (infoTable !== null) {\n    break;\n  }\n}\n\nif (infoTable === null) {\n  throw new Error('Could not find the \"Approbateurs\" header cell.');\n}\n\nconst values = infoTable.values;\n\n// The row right below the header holds the approver's name.\nconst nameRowIndex = headerRowIndex + 1;\n\n// Find the row (same column as \"Approbateurs\") holding \"Approuv\u00e9 le :\".\nlet approuveRowIndex = -1;\nfor (let r = nameRowIndex; r < values.length; r++) {\n  const cellText = values[r][approbateursColIndex] || \"\";\n  if (cellText.indexOf(\"Approuv\u00e9 le\") !== -1) {\n    approuveRowIndex = r;\n    break;\n  }\n}\n\nif (approuveRowIndex === -1) {\n  throw new Error('Could not find the \"Approuv\u00e9 le\" cell.');\n}\n\n// Add the approver's name to the (currently blank/space-only) cell below\n// the \"Approbateurs\" header.\nconst nameCell = infoTable.getCell(nameRowIndex, approbateursColIndex);\nnameCell.body.insertText(\"Fanny LAJEUNESSE\", Word.InsertLocation.end);\n\n// Add the approval date right after \"Approuv\u00e9 le :  \".\nconst approuveCell = infoTable.getCell(approuveRowIndex, approbateursColIndex);\napprouveCell.body.insertText(\"29/10\", Word.InsertLocation.end);\n\nawait context.sync();\n", "ps1": "# Fill in the approval info in the \"Auteurs / Approbateurs / Validation\"\n# table near the top of the document:\n#   - the \"Approbateurs\" value cell (currently just a single space) gets\n#     the approver's name \"Fanny LAJEUNESSE\" appended;\n#   - the \"Approuv\u00e9 le :  \" value cell gets the approval date \"29/10\"\n#     appended.\n\n$d = $word.ActiveDocument\n$tables = $d.Tables\n\n$foundTableIndex = -1\n$headerRow = -1\n$headerCol = -1\n\n# Locate the \"Approbateurs\" header cell (first match in document order) so\n# the edit is resilient to the table not being exactly where we expect it.\nfor ($ti = 1; $ti -le $tables.Count -and $foundTableIndex -eq -1; $ti++) {\n  $t = $tables.Item($ti)\n  $rows = $t.Rows.Count\n  $cols = $t.Columns.Count\n  for ($r = 1; $r -le $rows -and $foundTableIndex -eq -1; $r++) {\n    for ($c = 1; $c -le $cols -and $foundTableIndex -eq -1; $c++) {\n      $cellText = ($t.Cell($r, $c).Range.Text -replace \"[\\r\\a]\", \"\").Trim()\n      if ($cellText -eq \"Approbateurs\") {\n        $foundTableIndex = $ti\n        $headerRow = $r\n        $headerCol = $c\n      }\n    }\n  }\n}\n\nif ($foundTableIndex -eq -1) {\n  throw 'Could not find the \"Approbateurs\" header cell.'\n}\n\n$infoTable = $tables.Item($foundTableIndex)\n\n# The row right below the header holds the approver's name.\n$nameRow = $headerRow + 1\n$nameCell = $infoTable.Cell($nameRow, $headerCol)\n$nameCell.Range.InsertAfter(\"Fanny LAJEUNESSE\")\n\n# Find the row (same column as \"Approbateurs\") holding \"Approuv\u00e9 le :\".\n$approveRow = -1\n$rowCount = $infoTable.Rows.Count\nfor ($r = $nameRow; $r -le $rowCount -and $approveRow -eq -1; $r++) {\n  $txt = ($infoTable.Cell($r, $headerCol).Range.Text -replace \"[\\r\\a]\", \"\")\n  if ($txt.Contains(\"Approuv\u00e9 le\")) {\n    $approveRow = $r\n  }\n}\n\nif ($approveRow -eq -1) {\n  throw 'Could not find the \"Approuv\u00e9 le\" cell.'\n}\n\n# Add the approval date right after \"Approuv\u00e9 le :  \".\n$dateCell = $infoTable.Cell($approveRow, $headerCol)\n$dateCell.Range.InsertAfter(\"29/10\")\n"}
